$d = $word.ActiveDocument

# The three inline pictures in the headers/footers were re-ordered/renamed:
#   footer1.xml : image1.png -> image2.png  (Pearson logo, wp:docPr id=3 / pic:cNvPr id=0)
#   footer2.xml : image1.png -> image2.png  (Pearson logo, wp:docPr id=2 / pic:cNvPr id=0)
#   header1.xml : image2.jpg -> image1.jpg  (BTec logo,    wp:docPr id=1 / pic:cNvPr id=0)
#
# The Word object model does not expose a writable "Name" property on
# InlineShape that keeps both the <wp:docPr> AND the nested <pic:cNvPr>
# name attributes in sync (only <wp:docPr> gets updated that way), so the
# rename is applied directly against the package's canonical WordOpenXML,
# which covers every part (document + headers + footers) in one shot.

$xml = $d.WordOpenXML

$replacements = @(
    @{
        old = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/>'
        new = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/>'
    },
    @{
        old = '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
        new = '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
    },
    @{
        old = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>'
        new = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>'
    },
    @{
        old = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>'
        new = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>'
    },
    @{
        old = '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
        new = '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'
    }
)

foreach ($r in $replacements) {
    $xml = $xml.Replace($r.old, $r.new)
}

$d.WordOpenXML = $xml

Write-Output "done"
